$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 48, shifting existing rows 48-74 down to 49-75,
# carrying over formatting (date style in column D) from the row being pushed down.
$ws.Range("A48:R48").Insert()

# Populate the newly inserted row 48 with the new weekly price record.
$ws.Range("A48").Value = 5
$ws.Range("B48").Value = "Macroferia Regional de Talca"
$ws.Range("C48").Value = "Maule"
$ws.Range("D48").Value = 44176
$ws.Range("E48").Value = 7
$ws.Range("F48").Value = 100112030
$ws.Range("G48").Value = "Poroto granado"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 100
$ws.Range("K48").Value = 38000
$ws.Range("L48").Value = 38000
$ws.Range("M48").Value = 38000
$ws.Range("N48").Value = "$/saco 25 kilos"
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 1520
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"
